# Updated cryptos list on Sat Nov 23 09:25:40 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to treat the assigned value as literal
# text (matching the original inlineStr cells) without altering the
# cell's NumberFormat (stays 'General', as in the source file).

$ws.Range("D2").Value = "'" + '98.451.63'
$ws.Range("E2").Value = "'" + '  -0.19%  '
$ws.Range("D3").Value = "'" + '3.357.06'
$ws.Range("E3").Value = "'" + '  -0.71%  '
$ws.Range("E4").Value = "'" + '  +0.00%  '
$ws.Range("D5").Value = "'" + '257.76'
$ws.Range("E5").Value = "'" + '  -0.63%  '
$ws.Range("D6").Value = "'" + '667.37'
$ws.Range("E6").Value = "'" + '  +6.12%  '
$ws.Range("D7").Value = "'" + '1.54'
$ws.Range("E7").Value = "'" + '  +11.00%  '
$ws.Range("D8").Value = "'" + '0.462'
$ws.Range("E8").Value = "'" + '  +17.69%  '
$ws.Range("D9").Value = "'" + '1.10'
$ws.Range("E9").Value = "'" + '  +27.24%  '
$ws.Range("D11").Value = "'" + '3.351.92'
$ws.Range("E11").Value = "'" + '  -0.77%  '
$ws.Range("E12").Value = "'" + '  +5.69%  '
$ws.Range("D13").Value = "'" + '42.32'
$ws.Range("E13").Value = "'" + '  +16.52%  '
$ws.Range("D14").Value = "'" + '0.0000271'
$ws.Range("E14").Value = "'" + '  +8.72%  '
$ws.Range("D15").Value = "'" + '98.233.51'
$ws.Range("E15").Value = "'" + '  -0.32%  '
$ws.Range("D16").Value = "'" + '3.982.52'
$ws.Range("E16").Value = "'" + '  -0.56%  '
$ws.Range("D17").Value = "'" + '5.60'
$ws.Range("E17").Value = "'" + '  +1.88%  '
$ws.Range("D18").Value = "'" + '3.357.60'
$ws.Range("E18").Value = "'" + '  -0.14%  '
$ws.Range("D19").Value = "'" + '7.65'
$ws.Range("E19").Value = "'" + '  +25.76%  '
$ws.Range("D20").Value = "'" + '16.83'
$ws.Range("E20").Value = "'" + '  +10.32%  '
$ws.Range("B21").Value = "'" + 'BitcoinCash'
$ws.Range("C21").Value = "'" + 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = "'" + '535.22'
$ws.Range("E21").Value = "'" + '  +9.38%  '
$ws.Range("B22").Value = "'" + 'SuiNetwork'
$ws.Range("C22").Value = "'" + 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").Value = "'" + '3.57'
$ws.Range("E22").Value = "'" + '  -0.19%  '
$ws.Range("D23").Value = "'" + '10.45'
$ws.Range("E23").Value = "'" + '  +11.02%  '
$ws.Range("D24").Value = "'" + '0.444'
$ws.Range("E24").Value = "'" + '  +56.80%  '
$ws.Range("D25").Value = "'" + '0.0000214'
$ws.Range("E25").Value = "'" + '  +1.05%  '
$ws.Range("D26").Value = "'" + '102.70'
$ws.Range("E26").Value = "'" + '  +15.37%  '
$ws.Range("D27").Value = "'" + '6.24'
$ws.Range("E27").Value = "'" + '  +10.39%  '
$ws.Range("D28").Value = "'" + '12.64'
$ws.Range("E28").Value = "'" + '  +5.32%  '
$ws.Range("D29").Value = "'" + '3.545.67'
$ws.Range("E29").Value = "'" + '  -0.27%  '
$ws.Range("D30").Value = "'" + '0.151'
$ws.Range("E30").Value = "'" + '  +10.86%  '
$ws.Range("D31").Value = "'" + '0.999'
$ws.Range("E31").Value = "'" + '  -0.17%  '
$ws.Range("D32").Value = "'" + '11.10'
$ws.Range("E32").Value = "'" + '  +14.70%  '
$ws.Range("E33").Value = "'" + '  -0.90%  '
$ws.Range("E34").Value = "'" + '  +0.10%  '
$ws.Range("D35").Value = "'" + '29.50'
$ws.Range("E35").Value = "'" + '  +5.08%  '
$ws.Range("E36").Value = "'" + '  +17.23%  '
$ws.Range("D37").Value = "'" + '7.86'
$ws.Range("E37").Value = "'" + '  +7.23%  '
$ws.Range("D38").Value = "'" + '0.160'
$ws.Range("E38").Value = "'" + '  +6.26%  '
$ws.Range("D39").Value = "'" + '2.12'
$ws.Range("E39").Value = "'" + '  +7.83%  '
$ws.Range("D40").Value = "'" + '527.56'
$ws.Range("E40").Value = "'" + '  +4.95%  '
$ws.Range("B41").Value = "'" + 'Fetch.AI'
$ws.Range("C41").Value = "'" + 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").Value = "'" + '1.35'
$ws.Range("E41").Value = "'" + '  +5.97%  '
$ws.Range("B42").Value = "'" + 'WhiteBITCoin'
$ws.Range("C42").Value = "'" + 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").Value = "'" + '24.70'
$ws.Range("E42").Value = "'" + '  -0.84%  '
$ws.Range("B43").Value = "'" + 'VeChain'
$ws.Range("C43").Value = "'" + 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = "'" + '0.0442'
$ws.Range("E43").Value = "'" + '  +35.50%  '
$ws.Range("D44").Value = "'" + '3.73'
$ws.Range("E44").Value = "'" + '  -0.80%  '
$ws.Range("D45").Value = "'" + '0.832'
$ws.Range("E45").Value = "'" + '  +4.53%  '
$ws.Range("D46").Value = "'" + '3.42'
$ws.Range("E46").Value = "'" + '  +3.35%  '
$ws.Range("D48").Value = "'" + '8.02'
$ws.Range("E48").Value = "'" + '  +20.34%  '
$ws.Range("B49").Value = "'" + 'Stacks'
$ws.Range("C49").Value = "'" + 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = "'" + '2.07'
$ws.Range("E49").Value = "'" + '  +6.61%  '
$ws.Range("B50").Value = "'" + 'Filecoin'
$ws.Range("C50").Value = "'" + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = "'" + '5.17'
$ws.Range("E50").Value = "'" + '  +11.32%  '
$ws.Range("D51").Value = "'" + '1.53'
$ws.Range("E51").Value = "'" + '  +11.95%  '
